# Apply cryptos list update (prices/volumes refreshed; rows 40-41 and 48-49 swapped order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.607.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.007.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.002.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.488.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.667.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.001.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.672"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "457.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.191.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0792"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0381"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.68%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.242"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0488"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.15%  "
